# Updates "想去人数" (want-to-go count, column F) and, for one event whose
# ticket price changed, "最低票价" (min price, column G) on the "展览"
# sheet (sheet1) and the "全部类型" sheet (sheet4). Both sheets list the
# same events; "全部类型" simply has one extra leading row, so its row
# numbers for the same event are offset by +1 from row 12 onward.

$wb = $excel.ActiveWorkbook

# row -> new F ("想去人数") value, for the "展览" sheet
$sheet1F = @{
    6  = 173
    7  = 3726
    8  = 171
    12 = 74
    13 = 652
    14 = 153
    15 = 866
    16 = 53
    17 = 220
    18 = 147
    21 = 79
    22 = 3141
    23 = 5517
    26 = 504
    28 = 3179
    29 = 336
    30 = 2376
    32 = 506
    34 = 169
    35 = 235
    37 = 86
    38 = 488
    39 = 856
    41 = 24
    42 = 457
    44 = 525
}

# row -> new F ("想去人数") value, for the "全部类型" sheet
$sheet4F = @{
    6  = 173
    7  = 3726
    8  = 171
    13 = 74
    14 = 652
    15 = 153
    16 = 866
    17 = 53
    18 = 220
    19 = 147
    22 = 79
    23 = 3141
    24 = 5517
    27 = 504
    29 = 3179
    30 = 336
    31 = 2376
    33 = 506
    35 = 169
    36 = 235
    39 = 488
    40 = 856
    42 = 24
    43 = 457
    45 = 525
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1F.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1F[$row]
}
# Min ticket price for row 44 ("南昌·第四届龙年动漫展——暑假最后的狂欢") also changed
$ws1.Cells.Item(44, 7).Value = 55

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4F.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4F[$row]
}
# Same event, same min ticket price change, row shifted by +1 on this sheet
$ws4.Cells.Item(45, 7).Value = 55
